$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @(
    "Segmentname_FV2304",
    "Segmentgruppe_FV2304",
    "Segment_FV2304",
    "Datenelement_FV2304",
    "Segment ID_FV2304",
    "Code_FV2304",
    "Qualifier_FV2304",
    "Beschreibung_FV2304",
    "Bedingungsausdruck_FV2304",
    "Bedingung_FV2304",
    "diff",
    "Segmentname_FV2310",
    "Segmentgruppe_FV2310",
    "Segment_FV2310",
    "Datenelement_FV2310",
    "Segment ID_FV2310",
    "Code_FV2310",
    "Qualifier_FV2310",
    "Beschreibung_FV2310",
    "Bedingungsausdruck_FV2310",
    "Bedingung_FV2310"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Create the table (ListObject) over the data range, using existing headers.
$tableRange = $ws.Range("A1:U58")
$listObj = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$listObj.Name = "Table1"

# Freeze the top row (header row).
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
